$wb = $excel.ActiveWorkbook

# --- 1. Rename header cells on the existing sheets ---------------------
$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after "Monthly Trend" ----------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- 3. Header row ------------------------------------------------------
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Reuse the same header style already used on "Weekly Quantity" (bold,
# centered, bordered) by copy/paste-special-formats instead of rebuilding
# fonts/borders from scratch.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# --- 4. Data rows ---------------------------------------------------------
$forecastData = @(
    @(44969.99999999999, 95, 66.64517921673726, 122.6060505533512),
    @(45060.99999999999, 59, 32.02737128938729, 85.11237831988177),
    @(45088.99999999999, 47, 20.00873592528973, 75.26880257685001),
    @(45102.99999999999, 42, 14.57955397303586, 69.41715142840724),
    @(45116.99999999999, 36, 8.454478135009154, 63.4162998402088),
    @(45123.99999999999, 33, 5.739902344916689, 61.4222472343466),
    @(45130.99999999999, 31, 2.824371329270899, 58.52929255134833),
    @(45144.99999999999, 25, -0.5651089376789881, 52.75961566040882),
    @(45151.99999999999, 22, -3.231721713412995, 49.48768017785475),
    @(45158.99999999999, 19, -6.268896141869498, 45.59833900113),
    @(45165.99999999999, 17, -8.826077575089242, 43.37263550957393),
    @(45172.99999999999, 14, -13.68340438380553, 40.75088597845502),
    @(45179.99999999999, 11, -17.58308310710846, 38.80960081064077),
    @(45186.99999999999, 8, -18.43186449350769, 35.22336731408652),
    @(45193.99999999999, 5, -21.26654842947296, 32.06382438388501),
    @(45200.99999999999, 3, -25.30306563832194, 30.03475600009421)
)

for ($i = 0; $i -lt $forecastData.Length; $i++) {
    $r = $i + 2
    $row = $forecastData[$i]
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
}

# Reuse the date-column style ("s=2") used for the A-column dates on the
# other sheets for the new "ds" column.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A17").PasteSpecial(-4122)

$null = $wsForecast.Range("A1").Select()
